# BatchDepositSample.xlsx
# - Capitalize the data-field header labels on Sheet1 (row 1) so they
#   display nicely: code/issue/date/cms/type -> Code/Issue/Date/CMS/Type.
# - Update the remembered selection on Sheet1 from F5 to G9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "Code"
$ws.Range("B1").Value = "Issue"
$ws.Range("C1").Value = "Date"
$ws.Range("D1").Value = "CMS"
$ws.Range("E1").Value = "Type"

[void]$ws.Range("G9").Select()
